$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "62 x 54`v  5    4`v  ----`v6|    |`v2|    |"
$t.Cell(1, 2).Range.Text = "73 x 61`v  6    1`v  ----`v7|    |`v3|    |"
$t.Cell(1, 3).Range.Text = "39 x 67`v  6    7`v  ----`v3|    |`v9|    |"
$t.Cell(2, 1).Range.Text = "21 x 61`v  6    1`v  ----`v2|    |`v1|    |"
$t.Cell(2, 2).Range.Text = "96 x 24`v  2    4`v  ----`v9|    |`v6|    |"
$t.Cell(2, 3).Range.Text = "10 x 19`v  1    9`v  ----`v1|    |`v0|    |"
$t.Cell(3, 1).Range.Text = "99 x 71`v  7    1`v  ----`v9|    |`v9|    |"
$t.Cell(3, 2).Range.Text = "93 x 37`v  3    7`v  ----`v9|    |`v3|    |"
$t.Cell(3, 3).Range.Text = "83 x 79`v  7    9`v  ----`v8|    |`v3|    |"
$t.Cell(4, 1).Range.Text = "71 x 81`v  8    1`v  ----`v7|    |`v1|    |"
$t.Cell(4, 2).Range.Text = "75 x 94`v  9    4`v  ----`v7|    |`v5|    |"
$t.Cell(4, 3).Range.Text = "78 x 23`v  2    3`v  ----`v7|    |`v8|    |"
$t.Cell(5, 1).Range.Text = "59 x 86`v  8    6`v  ----`v5|    |`v9|    |"
$t.Cell(5, 2).Range.Text = "26 x 71`v  7    1`v  ----`v2|    |`v6|    |"
$t.Cell(5, 3).Range.Text = "27 x 92`v  9    2`v  ----`v2|    |`v7|    |"
